$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Add the new "Errors" worksheet after the existing sheets (Benchmarks, Notes)
# ---------------------------------------------------------------------------
$errorsSheet = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$errorsSheet.Name = "Errors"

# Title
$errorsSheet.Range("B2").Value = "GaussianBoxDelta"
$errorsSheet.Range("B2:B3").Font.Bold = $true
$errorsSheet.Range("B2:B3").Font.Size = 18

# Baseline delta row
$errorsSheet.Range("B4").Value = "Box"
$errorsSheet.Range("G4").Value = 0.0063895699999999998
$errorsSheet.Range("G4").NumberFormat = "0.0000"

# File label
$errorsSheet.Range("B6").Value = "Lucy.json"

# Table header
$errorsSheet.Range("B7").Value = "Rows"
$errorsSheet.Range("C7").Value = "Columns"
$errorsSheet.Range("E7").Value = "Passes"
$errorsSheet.Range("G7").Value = "MSE"
$errorsSheet.Range("G7").NumberFormat = "0.0000"

# Table data
$errorsSheet.Range("B8").Value = 512
$errorsSheet.Range("C8").Value = 1024

$passes = 1,2,3,4,5,10,20,100,10000
$mse = 0.61456118999999997943, 0.19017514999999998726, 0.04262644000000000155, 0.02035100000000000103, 0.01618034000000000139, 0.00678241999999999994, 0.00523959999999999964, 0.00547351999999999995, 0.00547351999999999995

for ($i = 0; $i -lt $passes.Length; $i++) {
    $row = 8 + $i
    $errorsSheet.Cells.Item($row, 5).Value = $passes[$i]
    $cell = $errorsSheet.Cells.Item($row, 7)
    $cell.Value = $mse[$i]
    $cell.NumberFormat = "0.0000"
}

$errorsSheet.Columns.Item(7).NumberFormat = "0.0000"

# ---------------------------------------------------------------------------
# Benchmarks sheet: update the saved selection / scroll position
# ---------------------------------------------------------------------------
$benchmarks = $wb.Worksheets.Item("Benchmarks")
$benchmarks.Activate()
$benchmarks.Range("A10").Select()

# ---------------------------------------------------------------------------
# Notes sheet: clear the saved selection
# ---------------------------------------------------------------------------
$notes = $wb.Worksheets.Item("Notes")
$notes.Activate()
$notes.Range("A1").Select()

$errorsSheet.Activate()
